$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 3752851.75
$ws.Range("C7").Value = -15.53486910106365
$ws.Range("D7").Value = 3338
$ws.Range("E7").Value = 3338
$ws.Range("F7").Value = 1124.2815308568
$ws.Range("G7").Value = 19.84028158698696
